# Update "想去人数" (wanted-to-go attendee count) figures in the
# "展览" (Exhibition) and "全部类型" (All Types) worksheets, column F.
#
# Source: gh-pages output regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (column F = "想去人数") ---
$wsExhibit = $wb.Worksheets.Item("展览")
$exhibitUpdates = @{
    2  = 6901
    4  = 447
    5  = 75
    6  = 15
    7  = 548
    8  = 115
    9  = 110
    11 = 8
    12 = 46
    13 = 188
    16 = 1803
    17 = 37
    18 = 3511
    19 = 23
    20 = 239
    21 = 23
    22 = 2140
    23 = 206
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# --- Sheet "全部类型" (column F = "想去人数") ---
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    2  = 6901
    4  = 447
    5  = 75
    6  = 15
    8  = 548
    9  = 115
    10 = 110
    12 = 8
    13 = 46
    14 = 188
    17 = 1803
    18 = 37
    19 = 3511
    20 = 23
    21 = 239
    22 = 23
    23 = 2140
    24 = 206
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}

$wb.Save()
